$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 45
$ws.Range("B2").Value = 129
$ws.Range("B3").Value = 159
$ws.Range("B4").Value = 171
$ws.Range("B5").Value = 214
$ws.Range("B6").Value = 228
$ws.Range("B7").Value = 279
